$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C to make room for "Variance".
# This shifts old C (Integral) -> D, old D (Time) -> E.
$ws.Columns.Item(3).Insert()

# Header for new column C
$ws.Range("C1").Value = "Variance"

# Fill C2:C11 with the variance formula (B^2), shared across the range.
$ws.Range("C2:C11").Formula = "=B2^2"

# Row 13 (Avg): C13 should now average the new Variance column.
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14 (STD): remove B14 and C14 (std of raw data / variance no longer wanted).
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Row 15: RMS = SQRT(average variance)
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Column C width, to match bestFit width in the diff (closest reachable value).
$ws.Columns.Item(3).ColumnWidth = 11.33

# Update selection to match diff (activeCell C13)
$ws.Range("C13").Select()
